# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ciudades")

# --- Insert a new data row for "Asturias" right after "Malaga" (row 19), ---
# --- pushing "Gran Canaria" ... "El Hierro" down by one row.              ---
$ws.Rows.Item(19).Insert()

$ws.Range("A19").Value = "Asturias"
$ws.Range("B19").Value = 1236
$ws.Range("C19").Value = 90
$ws.Range("D19").Value = 1091
$ws.Range("E19").Value = 55

# --- The old "Asturias" row (now pushed down to row 26, since the sheet ---
# --- had an existing "Asturias" row at 25 before the insert) is now a   ---
# --- duplicate; remove it so each province appears once.                ---
$ws.Rows.Item(26).Delete()

# --- Standalone update to the "Murcia" row's figures. ---
$ws.Range("C33").Value = 20
$ws.Range("D33").Value = 920
$ws.Range("E33").Value = 34

# --- Update the "last updated" timestamp string. ---
$ws.Range("A1").Value = "Datos actualizados a 30 de Marzo de 2020 a las 21:50"
